$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the D (Price) column cells we touch are written as plain text,
# matching the source workbook where these look like numbers/dates but are
# actually inline string values (e.g. "59.015.98", "521.95", "0.570").

# Row 2: ('Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '59.200.39', '  +2.20%  ') -> ('Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '59.015.98', '  +1.42%  ')
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.015.98"
$ws.Range("E2").Value = "  +1.42%  "

# Row 3: ('Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '2.633.25', '  +4.28%  ') -> ('Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '2.629.52', '  +3.49%  ')
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.629.52"
$ws.Range("E3").Value = "  +3.49%  "

# Row 4: ('TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.00', '  -0.13%  ') -> ('TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.00', '  +0.28%  ')
$ws.Range("E4").Value = "  +0.28%  "

# Row 5: ('BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '523.76', '  +4.60%  ') -> ('BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '521.95', '  +3.39%  ')
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "521.95"
$ws.Range("E5").Value = "  +3.39%  "

# Row 6: ('Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '146.94', '  +3.13%  ') -> ('Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '146.49', '  +1.61%  ')
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.49"
$ws.Range("E6").Value = "  +1.61%  "

# Row 7: ('USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '0.995', '  -0.45%  ') -> ('USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '0.997', '  +0.01%  ')
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  +0.01%  "

# Row 8: ('XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.571', '  +2.17%  ') -> ('XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.570', '  +1.09%  ')
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.570"
$ws.Range("E8").Value = "  +1.09%  "

# Row 9: ('LidoStakedEther', 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth', '2.658.24', '  +4.81%  ') -> ('LidoStakedEther', 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth', '2.648.99', '  +4.09%  ')
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.648.99"
$ws.Range("E9").Value = "  +4.09%  "

# Row 10: ('Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '6.33', '  +4.07%  ') -> ('Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '6.35', '  +3.64%  ')
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.35"
$ws.Range("E10").Value = "  +3.64%  "

# Row 11: ('Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.106', '  +4.97%  ') -> ('Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.105', '  +3.40%  ')
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.105"
$ws.Range("E11").Value = "  +3.40%  "

# Row 12: ('Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.339', '  +3.32%  ') -> ('Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.338', '  +2.28%  ')
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.338"
$ws.Range("E12").Value = "  +2.28%  "

# Row 13: ('TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.126', '  -1.12%  ') -> ('TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.126', '  -1.14%  ')
$ws.Range("E13").Value = "  -1.14%  "

# Row 14: ('WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '3.099.83', '  +4.14%  ') -> ('WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '3.099.20', '  +4.13%  ')
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.099.20"
$ws.Range("E14").Value = "  +4.13%  "

# Row 15: ('WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '59.129.12', '  +2.02%  ') -> ('WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '59.049.03', '  +1.54%  ')
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.049.03"
$ws.Range("E15").Value = "  +1.54%  "

# Row 16: ('Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '21.06', '  +2.65%  ') -> ('Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '21.02', '  +1.66%  ')
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.02"
$ws.Range("E16").Value = "  +1.66%  "

# Row 17: ('ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.0000138', '  +3.45%  ') -> ('ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.0000137', '  +1.92%  ')
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000137"
$ws.Range("E17").Value = "  +1.92%  "

# Row 18: ('WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '2.652.32', '  +4.59%  ') -> ('WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '2.645.41', '  +4.52%  ')
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.645.41"
$ws.Range("E18").Value = "  +4.52%  "

# Row 19: ('BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '348.68', '  +2.52%  ') -> ('BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '347.47', '  +1.41%  ')
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "347.47"
$ws.Range("E19").Value = "  +1.41%  "

# Row 20: ('Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '4.53', '  +1.26%  ') -> ('Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '4.52', '  +0.64%  ')
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.52"
$ws.Range("E20").Value = "  +0.64%  "

# Row 21: ('Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '10.34', '  +3.73%  ') -> ('Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '10.31', '  +2.25%  ')
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.31"
$ws.Range("E21").Value = "  +2.25%  "

# Row 22: ('Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '6.18', '  +5.09%  ') -> ('Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '6.18', '  +4.24%  ')
$ws.Range("E22").Value = "  +4.24%  "

# Row 23: ('Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '0.998', '  -0.17%  ') -> ('Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '0.997', '  -0.42%  ')
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.42%  "

# Row 24: ('Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '61.80', '  +2.70%  ') -> ('Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '61.76', '  +2.40%  ')
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.76"
$ws.Range("E24").Value = "  +2.40%  "

# Row 25: ('Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.419', '  +3.40%  ') -> ('Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.419', '  +2.36%  ')
$ws.Range("E25").Value = "  +2.36%  "

# Row 26: ('Kaspa', 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas', '0.164', '  +4.76%  ') -> ('Kaspa', 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas', '0.164', '  +4.07%  ')
$ws.Range("E26").Value = "  +4.07%  "

# Row 27: ('Binance-PegBSC-USD', 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd', '0.993', '  -0.57%  ') -> ('Binance-PegBSC-USD', 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd', '0.996', '  -0.19%  ')
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.996"
$ws.Range("E27").Value = "  -0.19%  "

# Row 28: ('PEPE', 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe', '0.0₃0812', '  +5.09%  ') -> ('PEPE', 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe', '0.0₃0807', '  +2.81%  ')
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0807"
$ws.Range("E28").Value = "  +2.81%  "

# Row 29: ('InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '7.14', '  +4.16%  ') -> ('InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '7.13', '  +3.00%  ')
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.13"
$ws.Range("E29").Value = "  +3.00%  "

# Row 30: ('USDe', 'https://coinranking.com/coin/exbfr2U-0+usde-usde', '0.998', '  -0.23%  ') -> ('USDe', 'https://coinranking.com/coin/exbfr2U-0+usde-usde', '0.998', '  -0.04%  ')
$ws.Range("E30").Value = "  -0.04%  "

# Row 31: ('Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '6.28', '  +7.77%  ') -> ('Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '6.27', '  +6.65%  ')
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.27"
$ws.Range("E31").Value = "  +6.65%  "

# Row 32: ('PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '1.59', '  +4.80%  ') -> ('PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '1.59', '  +3.89%  ')
$ws.Range("E32").Value = "  +3.89%  "

# Row 33: ('EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '19.01', '  +3.71%  ') -> ('EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '18.98', '  +2.96%  ')
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.98"
$ws.Range("E33").Value = "  +2.96%  "

# Row 34: ('Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '149.97', '  +0.68%  ') -> ('Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '149.95', '  +0.61%  ')
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.95"
$ws.Range("E34").Value = "  +0.61%  "

# Row 35: ('SuiNetwork', 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui', '0.979', '  +9.86%  ') -> ('SuiNetwork', 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui', '0.974', '  +7.85%  ')
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.974"
$ws.Range("E35").Value = "  +7.85%  "

# Row 36: ('NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '4.02', '  +4.65%  ') -> ('NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '4.01', '  +2.85%  ')
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.01"
$ws.Range("E36").Value = "  +2.85%  "

# Row 37: ('ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '1.14', '  +3.86%  ') -> ('ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '1.14', '  +2.99%  ')
$ws.Range("E37").Value = "  +2.99%  "

# Row 38: ('OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '36.82', '  +2.78%  ') -> ('OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '36.75', '  +2.69%  ')
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.75"
$ws.Range("E38").Value = "  +2.69%  "

# Row 39: ('Fetch.AI', 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet', '0.850', '  +5.69%  ') -> ('Fetch.AI', 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet', '0.849', '  +3.15%  ')
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.849"
$ws.Range("E39").Value = "  +3.15%  "

# Row 40: ('Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '3.70', '  +6.17%  ') -> ('Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '3.69', '  +5.07%  ')
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.69"
$ws.Range("E40").Value = "  +5.07%  "

# Row 41: ('Stacks', 'https://coinranking.com/coin/mMPrMcB7+stacks-stx', '1.42', '  +4.83%  ') -> ('Stacks', 'https://coinranking.com/coin/mMPrMcB7+stacks-stx', '1.43', '  +3.69%  ')
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.43"
$ws.Range("E41").Value = "  +3.69%  "

# Row 42: ('Bittensor', 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao', '279.80', '  +1.31%  ') -> ('Bittensor', 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao', '278.43', '  -0.88%  ')
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "278.43"
$ws.Range("E42").Value = "  -0.88%  "

# Row 43: ('Mantle', 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt', '0.611', '  +2.55%  ') -> ('FirstDigitalUSD', 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd', '0.996', '  -0.11%  ')
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.996"
$ws.Range("E43").Value = "  -0.11%  "

# Row 44: ('Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.0988', '  +0.35%  ') -> ('Mantle', 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt', '0.610', '  +2.31%  ')
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.610"
$ws.Range("E44").Value = "  +2.31%  "

# Row 45: ('FirstDigitalUSD', 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd', '0.994', '  -0.74%  ') -> ('Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.0986', '  -0.37%  ')
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0986"
$ws.Range("E45").Value = "  -0.37%  "

# Row 46: ('EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '19.64', '  +6.64%  ') -> ('EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '19.63', '  +4.71%  ')
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.63"
$ws.Range("E46").Value = "  +4.71%  "

# Row 47: ('Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.0526', '  -0.15%  ') -> ('Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.0526', '  -0.84%  ')
$ws.Range("E47").Value = "  -0.84%  "

# Row 48: ('RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '4.76', '  +6.83%  ') -> ('VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.0230', '  +1.95%  ')
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0230"
$ws.Range("E48").Value = "  +1.95%  "

# Row 49: ('VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.0230', '  +3.51%  ') -> ('WhiteBITCoin', 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt', '10.28', '  +0.41%  ')
$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.28"
$ws.Range("E49").Value = "  +0.41%  "

# Row 50: ('WhiteBITCoin', 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt', '10.30', '  +0.32%  ') -> ('RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '4.70', '  +3.94%  ')
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.70"
$ws.Range("E50").Value = "  +3.94%  "

# Row 51: ('Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '1.993.50', '  +5.75%  ') -> ('Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '1.992.04', '  +5.09%  ')
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.992.04"
$ws.Range("E51").Value = "  +5.09%  "
